$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows whose match data changed (re-ordering / refreshed odds) ---
# Row 2
$ws.Range("B2").Value = 6811743
$ws.Range("F2").Value = 'FC Thun'
$ws.Range("G2").Value = 'Stade Nyonnais'
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 'D'
$ws.Range("K2").Value = 1.615
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 1.65
$ws.Range("O2").Value = 4.333
$ws.Range("P2").Value = 4.5
$ws.Range("Q2").Value = -0.75
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 1.85
$ws.Range("V2").Value = 1.95
$ws.Range("W2").Value = -1
$ws.Range("X2").Value = 3.333
$ws.Range("Z2").Value = -1
$ws.Range("AA2").Value = 1.025
$ws.Range("AB2").Value = -1
$ws.Range("AC2").Value = 0.95

# Row 3
$ws.Range("B3").Value = 6811909
$ws.Range("F3").Value = 'Neuchatel Xamax'
$ws.Range("G3").Value = 'FC Schaffhausen'
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 'H'
$ws.Range("K3").Value = 2.4
$ws.Range("L3").Value = 3.5
$ws.Range("M3").Value = 2.7
$ws.Range("N3").Value = 2.05
$ws.Range("O3").Value = 3.75
$ws.Range("P3").Value = 3.5
$ws.Range("Q3").Value = -0.25
$ws.Range("T3").Value = 2.75
$ws.Range("U3").Value = 1.825
$ws.Range("V3").Value = 1.975
$ws.Range("W3").Value = 1.05
$ws.Range("X3").Value = -1
$ws.Range("Z3").Value = 0.7749999999999999
$ws.Range("AA3").Value = -1
$ws.Range("AB3").Value = 0.4125
$ws.Range("AC3").Value = -0.5

# Row 10
$ws.Range("B10").Value = 6811912
$ws.Range("F10").Value = 'Wil 1900'
$ws.Range("G10").Value = 'FC Vaduz'
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = 'H'
$ws.Range("K10").Value = 1.727
$ws.Range("L10").Value = 3.75
$ws.Range("M10").Value = 4
$ws.Range("N10").Value = 2.2
$ws.Range("O10").Value = 3.8
$ws.Range("P10").Value = 3
$ws.Range("Q10").Value = -0.25
$ws.Range("R10").Value = 1.9
$ws.Range("S10").Value = 1.9
$ws.Range("T10").Value = 3
$ws.Range("U10").Value = 1.975
$ws.Range("V10").Value = 1.825
$ws.Range("W10").Value = 1.2
$ws.Range("X10").Value = -1
$ws.Range("Z10").Value = 0.8999999999999999
$ws.Range("AA10").Value = -1
$ws.Range("AB10").Value = 0.9750000000000001
$ws.Range("AC10").Value = -1

# Row 11
$ws.Range("B11").Value = 6811424
$ws.Range("F11").Value = 'FC Baden'
$ws.Range("G11").Value = 'AC Bellinzona'
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 'D'
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 3.5
$ws.Range("M11").Value = 3.2
$ws.Range("N11").Value = 2.55
$ws.Range("O11").Value = 3.4
$ws.Range("P11").Value = 2.7
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 1.825
$ws.Range("S11").Value = 1.975
$ws.Range("T11").Value = 2.75
$ws.Range("U11").Value = 1.85
$ws.Range("V11").Value = 1.95
$ws.Range("W11").Value = -1
$ws.Range("X11").Value = 2.4
$ws.Range("Z11").Value = 0
$ws.Range("AA11").Value = 0
$ws.Range("AB11").Value = -1
$ws.Range("AC11").Value = 0.95

# Row 28
$ws.Range("B28").Value = 6811738
$ws.Range("F28").Value = 'Aarau'
$ws.Range("G28").Value = 'Stade Nyonnais'
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 2
$ws.Range("J28").Value = 'A'
$ws.Range("K28").Value = 1.727
$ws.Range("L28").Value = 3.5
$ws.Range("M28").Value = 4
$ws.Range("N28").Value = 1.55
$ws.Range("O28").Value = 4.5
$ws.Range("P28").Value = 5.5
$ws.Range("Q28").Value = -1
$ws.Range("T28").Value = 3.25
$ws.Range("U28").Value = 1.8
$ws.Range("V28").Value = 2
$ws.Range("X28").Value = -1
$ws.Range("Y28").Value = 4.5
$ws.Range("AC28").Value = 1

# Row 29
$ws.Range("B29").Value = 6811430
$ws.Range("F29").Value = 'FC Sion'
$ws.Range("G29").Value = 'FC Baden'
$ws.Range("H29").Value = 1
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = 'D'
$ws.Range("K29").Value = 1.3
$ws.Range("L29").Value = 4.5
$ws.Range("M29").Value = 8
$ws.Range("N29").Value = 1.222
$ws.Range("O29").Value = 6
$ws.Range("P29").Value = 13
$ws.Range("Q29").Value = -2
$ws.Range("T29").Value = 3.5
$ws.Range("U29").Value = 1.9
$ws.Range("V29").Value = 1.9
$ws.Range("X29").Value = 5
$ws.Range("Y29").Value = -1
$ws.Range("AC29").Value = 0.8999999999999999

# Row 45
$ws.Range("B45").Value = 6811735
$ws.Range("F45").Value = 'AC Bellinzona'
$ws.Range("G45").Value = 'Stade Nyonnais'
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 'H'
$ws.Range("K45").Value = 2.75
$ws.Range("L45").Value = 3.2
$ws.Range("M45").Value = 2.4
$ws.Range("O45").Value = 3.25
$ws.Range("P45").Value = 2.625
$ws.Range("R45").Value = 1.975
$ws.Range("S45").Value = 1.825
$ws.Range("U45").Value = 1.975
$ws.Range("V45").Value = 1.825
$ws.Range("W45").Value = 1.7
$ws.Range("X45").Value = -1
$ws.Range("Z45").Value = 0.9750000000000001
$ws.Range("AA45").Value = -1
$ws.Range("AC45").Value = 0.825

# Row 46
$ws.Range("B46").Value = 6811435
$ws.Range("F46").Value = 'FC Baden'
$ws.Range("G46").Value = 'FC Schaffhausen'
$ws.Range("I46").Value = 1
$ws.Range("J46").Value = 'D'
$ws.Range("K46").Value = 2.5
$ws.Range("L46").Value = 3.4
$ws.Range("M46").Value = 2.5
$ws.Range("O46").Value = 3.75
$ws.Range("P46").Value = 2.375
$ws.Range("R46").Value = 2.025
$ws.Range("S46").Value = 1.775
$ws.Range("U46").Value = 1.775
$ws.Range("V46").Value = 2.025
$ws.Range("W46").Value = -1
$ws.Range("X46").Value = 2.75
$ws.Range("Z46").Value = 0
$ws.Range("AA46").Value = 0
$ws.Range("AC46").Value = 1.025

# Row 57
$ws.Range("B57").Value = 6811732
$ws.Range("F57").Value = 'Stade Nyonnais'
$ws.Range("G57").Value = 'AC Bellinzona'
$ws.Range("H57").Value = 2
$ws.Range("I57").Value = 3
$ws.Range("K57").Value = 1.909
$ws.Range("L57").Value = 3.4
$ws.Range("M57").Value = 3.6
$ws.Range("N57").Value = 2
$ws.Range("O57").Value = 3.5
$ws.Range("P57").Value = 3.8
$ws.Range("Q57").Value = -0.5
$ws.Range("R57").Value = 2
$ws.Range("S57").Value = 1.8
$ws.Range("T57").Value = 2.5
$ws.Range("U57").Value = 1.8
$ws.Range("V57").Value = 2
$ws.Range("Y57").Value = 2.8
$ws.Range("AA57").Value = 0.8
$ws.Range("AB57").Value = 0.8
$ws.Range("AC57").Value = -1

# Row 58
$ws.Range("B58").Value = 6811439
$ws.Range("F58").Value = 'Wil 1900'
$ws.Range("G58").Value = 'FC Baden'
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 1
$ws.Range("K58").Value = 1.45
$ws.Range("L58").Value = 4.333
$ws.Range("M58").Value = 5.5
$ws.Range("N58").Value = 1.45
$ws.Range("O58").Value = 4.5
$ws.Range("P58").Value = 7
$ws.Range("Q58").Value = -1.25
$ws.Range("R58").Value = 1.9
$ws.Range("S58").Value = 1.9
$ws.Range("T58").Value = 3
$ws.Range("U58").Value = 1.925
$ws.Range("V58").Value = 1.875
$ws.Range("Y58").Value = 6
$ws.Range("AA58").Value = 0.8999999999999999
$ws.Range("AB58").Value = -1
$ws.Range("AC58").Value = 0.875

# Row 108
$ws.Range("B108").Value = 7617772
$ws.Range("F108").Value = 'FC Vaduz'
$ws.Range("G108").Value = 'FC Sion'
$ws.Range("I108").Value = 2
$ws.Range("J108").Value = 'A'
$ws.Range("K108").Value = 3.75
$ws.Range("L108").Value = 3.6
$ws.Range("M108").Value = 1.833
$ws.Range("N108").Value = 5.25
$ws.Range("O108").Value = 4
$ws.Range("P108").Value = 1.615
$ws.Range("Q108").Value = 0.75
$ws.Range("R108").Value = 2
$ws.Range("S108").Value = 1.8
$ws.Range("T108").Value = 2.75
$ws.Range("U108").Value = 1.825
$ws.Range("V108").Value = 1.975
$ws.Range("W108").Value = -1
$ws.Range("Y108").Value = 0.615
$ws.Range("Z108").Value = -0.5
$ws.Range("AA108").Value = 0.4
$ws.Range("AB108").Value = 0.4125
$ws.Range("AC108").Value = -0.5

# Row 109
$ws.Range("B109").Value = 7617773
$ws.Range("F109").Value = 'FC Thun'
$ws.Range("G109").Value = 'Aarau'
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 'H'
$ws.Range("K109").Value = 1.727
$ws.Range("L109").Value = 3.8
$ws.Range("M109").Value = 4
$ws.Range("N109").Value = 1.7
$ws.Range("O109").Value = 4.2
$ws.Range("P109").Value = 4.5
$ws.Range("Q109").Value = -0.75
$ws.Range("R109").Value = 1.85
$ws.Range("S109").Value = 1.95
$ws.Range("T109").Value = 3
$ws.Range("U109").Value = 1.9
$ws.Range("V109").Value = 1.9
$ws.Range("W109").Value = 0.7
$ws.Range("Y109").Value = -1
$ws.Range("Z109").Value = 0.425
$ws.Range("AA109").Value = -0.5
$ws.Range("AB109").Value = -1
$ws.Range("AC109").Value = 0.8999999999999999

# Row 139
$ws.Range("B139").Value = 7617795
$ws.Range("E139").Value = 45388.63541666666
$ws.Range("F139").Value = 'FC Sion'
$ws.Range("G139").Value = 'FC Vaduz'
$ws.Range("K139").Value = 1.5
$ws.Range("L139").Value = 4.2
$ws.Range("M139").Value = 5.25
$ws.Range("N139").Value = 1.444
$ws.Range("O139").Value = 4.5
$ws.Range("P139").Value = 7
$ws.Range("Q139").Value = -1.25
$ws.Range("R139").Value = 1.925
$ws.Range("S139").Value = 1.875
$ws.Range("T139").Value = 3
$ws.Range("U139").Value = 1.85
$ws.Range("V139").Value = 1.95

# --- Remove the trailing fixtures that were dropped from the feed ---
$ws.Range("A140:A142").EntireRow.Delete()

